$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Prelevante" name was re-entered as "Surname Name" instead of "Name Surname"
$ws.Cells.Item(11, 2).Value = "D'alesio Katia"

# New withdrawal (28 aprile / same date serial used already in row 11)
$ws.Cells.Item(12, 1).Value = 43217
$ws.Cells.Item(12, 1).NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(12, 1).VerticalAlignment = -4160

$ws.Cells.Item(12, 2).Value = "Lanzi Patrizia"
$ws.Cells.Item(12, 2).NumberFormat = "@"
$ws.Cells.Item(12, 2).VerticalAlignment = -4160

$ws.Cells.Item(12, 3).Value = "Bindello Nero"
$ws.Cells.Item(12, 3).NumberFormat = "@"
$ws.Cells.Item(12, 3).VerticalAlignment = -4160

$ws.Cells.Item(12, 4).Value = "N°."
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).VerticalAlignment = -4160

$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 5).VerticalAlignment = -4160
